# Update "想去人数" (attendance) figures on the "展览" and "全部类型" sheets
# F2: 9417 -> 9433
# F4: 23   -> 24
# F5: 517  -> 520

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9433
    $ws.Range("F4").Value = 24
    $ws.Range("F5").Value = 520
}
